# games/dqm2/translation-progress.xlsx - update trans progress xls
#
# Field 61-65 (row 12) was already "collected"; this edit:
#   - moves the "collected" status text into C12 (reusing the standard
#     "translation collected" label used elsewhere in the sheet)
#   - frees up the old "62-65 collected" label text and repurposes it as a
#     new field-range heading, "Field 66 - 70"
#   - adds a new row 13 for that field range, mirroring row 12's A/B layout
#     (same translator id, still pending collection so C13 stays blank)
#   - gives the blank C11/C13 cells matching "blank status cell" formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C12: "62-65已收回" -> becomes the normal "已收回" status used elsewhere
# (copy C8's look, which already carries that exact status text/format)
$ws.Range("C8").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = "翻译已收回"

# The shared string previously shown in C12 ("62-65已收回") is now unused;
# repurpose it in place as the new field-range header for row 13.
$ws.Range("A13").Value = "Field 66 - 70"

# Row 13 mirrors row 12's A/B formatting (field header + translator id).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = 9527

# C13 stays blank, styled like the other blank status cells (row 12's old
# look, now freed up) and C11 picks up the same "blank status" styling
# (matching C7, the other pending row directly above a collected one).
$ws.Range("C7").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C13").PasteSpecial(-4122)
